$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record columns (AD:AF), formatted like the
# other header cells (bold, centered, thin border) by copying AC1's style.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's record (Wins/Losses/Ties) for every player row.
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 53
    $ws.Range("AE$r").Value = 62
    $ws.Range("AF$r").Value = 0
}
